$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9480683207511902
$ws.Range("B1").Value = 3.091810703277588
$ws.Range("C1").Value = 2.756330966949463
$ws.Range("D1").Value = 1.566226601600647
$ws.Range("E1").Value = 1.203850746154785
